{"js": "// Commit: \"add a line in python\"\n// The document gains three new paragraphs at the very end of the body:\n//   1. an empty paragraph\n//   2. an empty paragraph\n//   3. a paragraph containing the text \"The note of python\"\n// All pre-existing content is left untouched.\n\nconst body = context.document.body;\n\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"The note of python\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Commit: \"add a line in python\"\n# The document gains three new paragraphs at the very end of the body:\n#   1. an empty paragraph\n#   2. an empty paragraph\n#   3. a paragraph containing the text \"The note of python\"\n# All pre-existing content is left untouched.\n\n$d = $word.ActiveDocument\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.InsertBefore(\"The note of python\")\n"}
